# 自动更新价格数据: insert today's price row at the top of the table,
# pushing the existing history down by one row (same prices, new date).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, copying formatting from the row below
# (the plain data rows) instead of the header row.
$ws.Rows("2:2").Insert(-4121, 1)

# New row 2: latest date, same commodity prices as the rest of the table.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-11-23"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# Make sure the inserted row carries no special formatting (matches the
# plain, unstyled data rows already present in the sheet).
$ws.Range("A2:D2").ClearFormats()
